$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Increased regen value on skip (I5): 7 -> 9, formulas M5/N5/O5 recalc automatically
$ws.Range("I5").Value = 9

# Move the class list (HOLY PRIEST / DISC PRIEST / HOLY PALADIN) up to rows 10-12
$ws.Range("A10").Value = "HOLY PRIEST"
$ws.Range("A11").Value = "DISC PRIEST"
$ws.Range("Q11").Value = "dispel"
$ws.Range("S11").Value = "mana burn"
$ws.Range("A12").Value = "HOLY PALADIN"

# Clear old positions of that list, and repurpose rows 17-21
$ws.Range("A17").ClearContents()
$ws.Range("A18").Value = "COMPS"
$ws.Range("A19").ClearContents()
$ws.Range("A20").Value = "BEAST HUNTER"
$ws.Range("A21").Value = "DEMON LOCK"

# New comps / notes section
$ws.Range("A22").Value = "RETRI PALADIN"

$ws.Range("J23").Value = "blessing of protection -> didn’t work when target was stunned"

$ws.Range("A24").Value = "DISC PRIEST"
$ws.Range("J24").Value = "viper sting reduces below 0"

$ws.Range("A25").Value = "AFFLI LOCK"
$ws.Range("J25").Value = "make rogues unable to skip?"

$ws.Range("A26").Value = "FROST MAGE"

$ws.Range("A28").Value = "SURV HUNTER"
$ws.Range("C28").Value = "no dmg"

$ws.Range("A29").Value = "DISC PRIEST"
$ws.Range("C29").Value = "lost vs all"

$ws.Range("A30").Value = "HOLY PALADIN"

$ws.Range("A32").Value = "SUB ROGUE"
$ws.Range("A33").Value = "FROST MAGE"
$ws.Range("A34").Value = "RETRI PALADIN"

$ws.Range("Q12").Select()
